# The deck's Slide Master and every Custom Layout carry a "Date
# Placeholder" shape whose auto-updating date field (type
# "datetimeFigureOut") was displaying a stale cached date
# ("2018-04-06"). Re-touching/re-saving the file in PowerPoint causes
# that cached date text to refresh to the new "as of" date
# ("2018-04-13"). Reproduce that by updating the date placeholder's
# text on the Slide Master and on every Custom Layout.

$p = $ppt.ActivePresentation

$newDate = "2018-04-13"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isPh = $false
        try { $isPh = $shp.HasTextFrame -and ($shp.PlaceholderFormat.Type -eq 16) } catch { $isPh = $false }
        if ($isPh) {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

# Slide Master.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes $newDate

# Every Custom Layout hanging off the master.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes $newDate
}
